$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.735929846763611
$ws.Range("B1").Value = 1.813096523284912
$ws.Range("C1").Value = 4.938595771789551
$ws.Range("D1").Value = 1.209845066070557
$ws.Range("E1").Value = 0.637754499912262
